{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the four paragraphs that make up the \"broken tutorials\" block:\n//   - \".../v_ISET.m -- BROKEN!\"        -> text replaced with \"\u00df\"\n//   - \".../v_isetcam.m -- BROKEN!\"     -> paragraph removed\n//   - \".../oi/v_icam_oi.m -- BROKEN!\"  -> paragraph removed\n//   - (empty paragraph right after)    -> paragraph removed\nconst items = paragraphs.items;\nlet targetIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"v_ISET.m -- BROKEN!\") !== -1) {\n    targetIndex = i;\n    break;\n  }\n}\n\nif (targetIndex === -1) {\n  throw new Error(\"Could not find the 'v_ISET.m -- BROKEN!' paragraph\");\n}\n\n// Replace the text of the first paragraph with \"\u00df\".\nitems[targetIndex].insertText(\"\u00df\", \"Replace\");\n\n// Delete the following three paragraphs (v_isetcam.m, oi/v_icam_oi.m, blank).\nitems[targetIndex + 1].delete();\nitems[targetIndex + 2].delete();\nitems[targetIndex + 3].delete();\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Find the paragraph that still reads \".../v_ISET.m -- BROKEN!\" so we don't\n# depend on a hard-coded paragraph index.\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*v_ISET.m -- BROKEN!*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find the 'v_ISET.m -- BROKEN!' paragraph\"\n}\n\n# Replace that paragraph's text with \"\u00df\" (keeps the run/paragraph formatting).\n$d.Paragraphs.Item($targetIndex).Range.Text = \"\u00df\"\n\n# Remove the next three paragraphs in one shot:\n#   - \".../v_isetcam.m -- BROKEN!\"\n#   - \".../oi/v_icam_oi.m -- BROKEN!\"\n#   - the blank paragraph that followed them\n$rangeStart = $d.Paragraphs.Item($targetIndex + 1).Range.Start\n$rangeEnd = $d.Paragraphs.Item($targetIndex + 3).Range.End\n$removeRange = $d.Range($rangeStart, $rangeEnd)\n$removeRange.Delete()\n"}
